# Slide 7 ("TRACTION & TRUST") edit:
# 3 AI-evaluator cards -> 4 AI-evaluator cards in a 2x2 grid, plus text
# updates, plus the gold highlighted-quote strip moving down to make room
# for the new second row of cards.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# ---------------------------------------------------------------------------
# 0. Grab handles to the trailing "extra" shapes (gold quote strip + footer)
#    up front, and rename the footer shapes out of the way FIRST so that the
#    renames we do later ("...18"->"...22" etc.) don't collide with names
#    that are still in use.
# ---------------------------------------------------------------------------
$goldBarOld   = $s.Shapes.Item("Rectangle 17")
$goldTextOld  = $s.Shapes.Item("TextBox 18")
$footerDiv    = $s.Shapes.Item("Rectangle 19")
$footerLeft   = $s.Shapes.Item("TextBox 20")
$footerMid    = $s.Shapes.Item("TextBox 21")
$footerPage   = $s.Shapes.Item("TextBox 22")

$footerDiv.Name  = "Rectangle 23"
$footerLeft.Name = "TextBox 24"
$footerMid.Name  = "TextBox 25"
$footerPage.Name = "TextBox 26"

# ---------------------------------------------------------------------------
# 1. Header title: "3社" -> "4社"
# ---------------------------------------------------------------------------
$s.Shapes.Item("TextBox 2").TextFrame.TextRange.Text = "主要AI 4社と文化機関が独立に評価——社会的信頼が構築されている"

# ---------------------------------------------------------------------------
# 2. Card 1 "Claude / Anthropic" (row 1, col 1) - reflow to wider/shorter card
# ---------------------------------------------------------------------------
$bg1 = $s.Shapes.Item("Rectangle 4")
$bg1.Left = 36 ; $bg1.Top = 82.8 ; $bg1.Width = 309.6 ; $bg1.Height = 97.2

$name1 = $s.Shapes.Item("TextBox 5")
$name1.Left = 46.8 ; $name1.Top = 88.56 ; $name1.Width = 288 ; $name1.Height = 15.84

$title1 = $s.Shapes.Item("TextBox 6")
$title1.Left = 46.8 ; $title1.Top = 105.84 ; $title1.Width = 288 ; $title1.Height = 18

$quote1 = $s.Shapes.Item("TextBox 7")
$quote1.Left = 46.8 ; $quote1.Top = 126 ; $quote1.Width = 288 ; $quote1.Height = 46.8

# ---------------------------------------------------------------------------
# 3. Card 2 "Gemini / Google" (row 1, col 2) - reflow + new headline/quote
# ---------------------------------------------------------------------------
$bg2 = $s.Shapes.Item("Rectangle 8")
$bg2.Left = 363.6 ; $bg2.Top = 82.8 ; $bg2.Width = 309.6 ; $bg2.Height = 97.2

$name2 = $s.Shapes.Item("TextBox 9")
$name2.Left = 374.4 ; $name2.Top = 88.56 ; $name2.Width = 288 ; $name2.Height = 15.84

$title2 = $s.Shapes.Item("TextBox 10")
$title2.Left = 374.4 ; $title2.Top = 105.84 ; $title2.Width = 288 ; $title2.Height = 18
$title2.TextFrame.TextRange.Text = "唯一無二の選択肢"

$quote2 = $s.Shapes.Item("TextBox 11")
$quote2.Left = 374.4 ; $quote2.Top = 126 ; $quote2.Width = 288 ; $quote2.Height = 46.8
$quote2.TextFrame.TextRange.Text = "「『永遠』を夢物語ではなく、現実的な予算と確かなロジックで構築した」"

# ---------------------------------------------------------------------------
# 4. Card 3 "ChatGPT / OpenAI" moves from (row 1, col 3) to (row 2, col 1)
#    + new headline/quote
# ---------------------------------------------------------------------------
$bg3 = $s.Shapes.Item("Rectangle 12")
$bg3.Left = 36 ; $bg3.Top = 188.64 ; $bg3.Width = 309.6 ; $bg3.Height = 97.2

$name3 = $s.Shapes.Item("TextBox 13")
$name3.Left = 46.8 ; $name3.Top = 194.4 ; $name3.Width = 288 ; $name3.Height = 15.84

$title3 = $s.Shapes.Item("TextBox 14")
$title3.Left = 46.8 ; $title3.Top = 211.68 ; $title3.Width = 288 ; $title3.Height = 18
$title3.TextFrame.TextRange.Text = "新基準として推奨"

$quote3 = $s.Shapes.Item("TextBox 15")
$quote3.Left = 46.8 ; $quote3.Top = 231.84 ; $quote3.Width = 288 ; $quote3.Height = 46.8
$quote3.TextFrame.TextRange.Text = "「技術的基盤・社会的意義・コストパフォーマンスのいずれも高く評価できる」"

# ---------------------------------------------------------------------------
# 5. NEW Card 4 "Grok / xAI" (row 2, col 2)
#    - the existing "Rectangle 16" shape (old gold quote-box background) is
#      repurposed as the card-4 background: it moves and gets the white
#      card fill / light border instead of the gold fill
#    - a fresh gold quote-box background is duplicated from it FIRST, so we
#      keep that gold styling available for step 6 below
#    - the 3 new text boxes are duplicated from card 2's (Gemini) text boxes
#      so they inherit the exact same run/paragraph formatting
# ---------------------------------------------------------------------------
$goldBoxNew = $s.Shapes.Item("Rectangle 16").Duplicate()
$goldBoxNew.Name = "Rectangle 20"

$bg4 = $s.Shapes.Item("Rectangle 16")
$bg4.Left = 363.6 ; $bg4.Top = 188.64 ; $bg4.Width = 309.6 ; $bg4.Height = 97.2
$bg4.Fill.ForeColor.RGB = 16777215    # FFFFFF (white), stored as 0x00BBGGRR
$bg4.Line.Visible = $true
$bg4.Line.Weight = 0.75
$bg4.Line.ForeColor.RGB = 15788258    # E2E8F0, stored as 0x00BBGGRR

$name4 = $s.Shapes.Item("TextBox 9").Duplicate()
$name4.Name = "TextBox 17"
$name4.Left = 374.4 ; $name4.Top = 194.4 ; $name4.Width = 288 ; $name4.Height = 15.84
$name4.TextFrame.TextRange.Text = "Grok / xAI"

$title4 = $s.Shapes.Item("TextBox 10").Duplicate()
$title4.Name = "TextBox 18"
$title4.Left = 374.4 ; $title4.Top = 211.68 ; $title4.Width = 288 ; $title4.Height = 18
$title4.TextFrame.TextRange.Text = "上位0.001%級"

$quote4 = $s.Shapes.Item("TextBox 11").Duplicate()
$quote4.Name = "TextBox 19"
$quote4.Left = 374.4 ; $quote4.Top = 231.84 ; $quote4.Width = 288 ; $quote4.Height = 46.8
$quote4.TextFrame.TextRange.Text = "「『本気で人類の記憶のあり方を変えるかもしれない』レベル」"

# ---------------------------------------------------------------------------
# 6. Gold highlighted-quote strip moves down (to below the new row 2 of
#    cards) and shrinks slightly; reuse the gold background duplicated in
#    step 5, rename+reposition the accent bar and quote text (handles
#    captured in step 0, before any renames happened).
# ---------------------------------------------------------------------------
$goldBg = $goldBoxNew
$goldBg.Left = 36 ; $goldBg.Top = 301.68 ; $goldBg.Width = 619.2 ; $goldBg.Height = 43.2

$goldBarOld.Name = "Rectangle 21"
$goldBarOld.Left = 36 ; $goldBarOld.Top = 301.68 ; $goldBarOld.Width = 4.32 ; $goldBarOld.Height = 43.2

$goldTextOld.Name = "TextBox 22"
$goldTextOld.Left = 54 ; $goldTextOld.Top = 307.44 ; $goldTextOld.Width = 594 ; $goldTextOld.Height = 31.68
